# Add analyser to cash
# Adds an "ID" tag (column A) to rows belonging to the cash / monetary-funds
# block (T002) and extends the existing T001 tag (already present on row
# 115) down through the rest of the "cash flow portrait" block (rows
# 116-124). Also fixes a typo in the E7 comment and adds two new
# annotation notes in column A (row 7 and row 8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New analyser id "T002" tagging the 货币资金 (cash) and 有息负债
#     (interest-bearing debt) blocks ---
$ws.Range("A4").Value = "T002"
$ws.Range("A5").Value = "T002"
$ws.Range("A6").Value = "T002"

# 有息负债 (interest-bearing debt) block also tagged T002.
$ws.Range("A48").Value = "T002"
$ws.Range("A49").Value = "T002"

# Row 7: fix the typo in the existing comment (流动资严重 -> 流动资金严重).
$ws.Range("E7").Value = "定期存款很多，其它货币资金很多，流动资金严重缺乏"

# Row 8: new note flagging that this line has no data in "ts".
$ws.Range("A8").Value = "ts没有此数据"

# Row 7: tack on a clarifying question in column A.
$ws.Range("A7").Value = "流动资金指什么"

# --- Extend the existing T001 tag (already on A115) down through the
#     rest of the cash-flow-portrait rows ---
$ws.Range("A116").Value = "T001"
$ws.Range("A117").Value = "T001"
$ws.Range("A118").Value = "T001"
$ws.Range("A119").Value = "T001"
$ws.Range("A120").Value = "T001"
$ws.Range("A121").Value = "T001"
$ws.Range("A122").Value = "T001"
$ws.Range("A123").Value = "T001"
$ws.Range("A124").Value = "T001"

# Leave the view roughly where the author ended up editing.
$ws.Range("A49").Select()
